$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 421.71054
$ws.Range("J17").Value = 421.71054
$ws.Range("L17").Value = 1265.13162
$ws.Range("N17").Value = -1601.13162
$ws.Range("H18").Value = 4333.3335
$ws.Range("I18").Value = 4000
$ws.Range("J18").Value = 4500
$ws.Range("K18").Value = 4000
$ws.Range("L18").Value = 4500
$ws.Range("M18").Value = -3716
$ws.Range("N18").Value = -5068
$ws.Range("H19").Value = 1258.4615
$ws.Range("I19").Value = 1175.6
$ws.Range("J19").Value = 1310.25
$ws.Range("K19").Value = 1175.6
$ws.Range("L19").Value = 1310.25
$ws.Range("M19").Value = -1000.6
$ws.Range("N19").Value = -1660.25
$ws.Range("H31").Value = 201.5
$ws.Range("I31").Value = 201.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 604.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -374.5
$ws.Range("N31").ClearContents()
$ws.Range("H32").Value = 1183.6666
$ws.Range("I32").Value = 900.5
$ws.Range("J32").Value = 1750
$ws.Range("K32").Value = 900.5
$ws.Range("L32").Value = 1750
$ws.Range("M32").Value = -574.5
$ws.Range("N32").Value = -2402
$ws.Range("H34").Value = 4128.875
$ws.Range("I34").Value = 4128.875
$ws.Range("K34").Value = 4128.875
$ws.Range("M34").Value = -3925.875
$ws.Range("H36").Value = 4128.875
$ws.Range("I36").Value = 4128.875
$ws.Range("K36").Value = 4128.875
$ws.Range("M36").Value = -3413.875
$ws.Range("H38").Value = 429.6875
$ws.Range("I38").Value = 220.22223
$ws.Range("K38").Value = 660.66669
$ws.Range("M38").Value = -288.66669
$ws.Range("H39").Value = 550
$ws.Range("I39").Value = 550
$ws.Range("K39").Value = 1650
$ws.Range("M39").Value = -1354
$ws.Range("H40").Value = 3499.75
$ws.Range("J40").Value = 4333.3335
$ws.Range("L40").Value = 4333.3335
$ws.Range("N40").Value = -4683.3335
$ws.Range("H41").Value = 485.2857
$ws.Range("I41").Value = 324.25
$ws.Range("K41").Value = 324.25
$ws.Range("M41").Value = 115.75
$ws.Range("H42").Value = 529.8182
$ws.Range("I42").Value = 181.33333
$ws.Range("J42").Value = 948
$ws.Range("K42").Value = 543.99999
$ws.Range("L42").Value = 2844
$ws.Range("M42").Value = -313.99999
$ws.Range("N42").Value = -3304
$ws.Range("H43").Value = 9981.799999999999
$ws.Range("J43").Value = 9981.25
$ws.Range("L43").Value = 9981.25
$ws.Range("N43").Value = -10119.25
$ws.Range("H51").Value = 7115.8125
$ws.Range("J51").Value = 7543.077
$ws.Range("L51").Value = 7543.077
$ws.Range("N51").Value = -8511.077000000001
$ws.Range("H52").Value = 916.6667
$ws.Range("J52").Value = 916.6667
$ws.Range("L52").Value = 2750.0001
$ws.Range("N52").Value = -3070.0001
$ws.Range("H53").Value = 1033.8334
$ws.Range("I53").Value = 764.6923
$ws.Range("J53").Value = 1733.6
$ws.Range("K53").Value = 764.6923
$ws.Range("L53").Value = 1733.6
$ws.Range("M53").Value = -127.6923
$ws.Range("N53").Value = -3007.6
$ws.Range("H54").Value = 7325
$ws.Range("I54").Value = 7325
$ws.Range("K54").Value = 7325
$ws.Range("M54").Value = -6839
$ws.Range("H55").Value = 239.33333
$ws.Range("I55").Value = 201.6
$ws.Range("J55").Value = 258.2
$ws.Range("K55").Value = 201.6
$ws.Range("L55").Value = 258.2
$ws.Range("M55").Value = 12.40000000000001
$ws.Range("N55").Value = -686.2
$ws.Range("H58").Value = 9999.5
$ws.Range("J58").Value = 9999.5
$ws.Range("L58").Value = 29998.5
$ws.Range("N58").Value = -30298.5
$ws.Range("H70").Value = 166586.23
$ws.Range("I70").Value = 233774.92
$ws.Range("J70").Value = 5333.4
$ws.Range("K70").Value = 701324.76
$ws.Range("L70").Value = 16000.2
$ws.Range("M70").Value = -701054.76
$ws.Range("N70").Value = -16540.2
$ws.Range("H73").Value = 166586.23
$ws.Range("I73").Value = 233774.92
$ws.Range("J73").Value = 5333.4
$ws.Range("K73").Value = 701324.76
$ws.Range("L73").Value = 16000.2
$ws.Range("M73").Value = -700388.76
$ws.Range("N73").Value = -17872.2
$ws.Range("H74").Value = 6832.8335
$ws.Range("J74").Value = 6499.25
$ws.Range("L74").Value = 6499.25
$ws.Range("N74").Value = -8371.25
$ws.Range("H77").Value = 6832.8335
$ws.Range("J77").Value = 6499.25
$ws.Range("L77").Value = 32496.25
$ws.Range("N77").Value = -41856.25
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H88").Value = 1001577
$ws.Range("I88").Value = 1999999
$ws.Range("J88").Value = 3155
$ws.Range("K88").Value = 1999999
$ws.Range("L88").Value = 3155
$ws.Range("M88").Value = -1999593
$ws.Range("N88").Value = -3967
$ws.Range("H91").Value = 1001577
$ws.Range("I91").Value = 1999999
$ws.Range("J91").Value = 3155
$ws.Range("K91").Value = 1999999
$ws.Range("L91").Value = 3155
$ws.Range("M91").Value = -1998595
$ws.Range("N91").Value = -5963
$ws.Range("H109").Value = 92665.336
$ws.Range("J109").Value = 92665.336
$ws.Range("L109").Value = 92665.336
$ws.Range("N109").Value = -95439.336
$ws.Range("H116").Value = 7148.8335
$ws.Range("I116").Value = 6978.6
$ws.Range("K116").Value = 6978.6
$ws.Range("M116").Value = -3536.6
$ws.Range("H118").Value = 303.2
$ws.Range("I118").Value = 303.2
$ws.Range("K118").Value = 909.5999999999999
$ws.Range("M118").Value = 747.4000000000001
$ws.Range("H126").Value = 67494.336
$ws.Range("J126").Value = 67494.336
$ws.Range("L126").Value = 67494.336
$ws.Range("N126").Value = -77374.336
$ws.Range("H127").Value = 1712.9231
$ws.Range("I127").Value = 1712.9231
$ws.Range("K127").Value = 5138.7693
$ws.Range("M127").Value = -178.7692999999999
$ws.Range("H132").Value = 3875.2778
$ws.Range("I132").Value = 3250.6667
$ws.Range("K132").Value = 9752.000100000001
$ws.Range("M132").Value = -7222.000100000001
$ws.Range("H135").Value = 8170.643
$ws.Range("I135").Value = 1106.9231
$ws.Range("K135").Value = 9962.3079
$ws.Range("M135").Value = -7427.3079
$ws.Range("H137").Value = 3230.75
$ws.Range("I137").Value = 1266.8334
$ws.Range("K137").Value = 3800.5002
$ws.Range("M137").Value = -1250.5002
$ws.Range("H138").Value = 2338.42
$ws.Range("I138").Value = 2860.7273
$ws.Range("J138").Value = 2191.1025
$ws.Range("K138").Value = 8582.1819
$ws.Range("L138").Value = 6573.3075
$ws.Range("M138").Value = -3442.1819
$ws.Range("N138").Value = -16853.3075

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1468.7142
$ws.Range("I45").Value = 1231.5
$ws.Range("J45").Value = 1943.1428
$ws.Range("K45").Value = 1231.5
$ws.Range("L45").Value = 1943.1428
$ws.Range("M45").Value = -854.5
$ws.Range("N45").Value = -2697.1428
$ws.Range("H74").Value = 1994.6086
$ws.Range("J74").Value = 2510.625
$ws.Range("L74").Value = 2510.625
$ws.Range("N74").Value = -4258.625
$ws.Range("H77").Value = 1994.6086
$ws.Range("J77").Value = 2510.625
$ws.Range("L77").Value = 12553.125
$ws.Range("N77").Value = -21289.125
$ws.Range("H88").Value = 1903.7273
$ws.Range("I88").Value = 2062.8333
$ws.Range("K88").Value = 2062.8333
$ws.Range("M88").Value = -1656.8333
$ws.Range("H91").Value = 1903.7273
$ws.Range("I91").Value = 2062.8333
$ws.Range("K91").Value = 2062.8333
$ws.Range("M91").Value = -658.8332999999998
$ws.Range("H132").Value = 3475.8572
$ws.Range("I132").Value = 3074.111
$ws.Range("K132").Value = 9222.332999999999
$ws.Range("M132").Value = -6692.332999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1132.2222
$ws.Range("I20").Value = 1149
$ws.Range("K20").Value = 1149
$ws.Range("M20").Value = -902
$ws.Range("H22").Value = 499.66666
$ws.Range("I22").Value = 498.5
$ws.Range("J22").Value = 502
$ws.Range("K22").Value = 498.5
$ws.Range("L22").Value = 502
$ws.Range("M22").Value = -325.5
$ws.Range("N22").Value = -848
$ws.Range("H64").Value = 999
$ws.Range("I64").Value = 999
$ws.Range("K64").Value = 999
$ws.Range("M64").Value = -774
$ws.Range("H67").Value = 999
$ws.Range("I67").Value = 999
$ws.Range("K67").Value = 999
$ws.Range("M67").Value = -219
$ws.Range("H86").Value = 2867.25
$ws.Range("I86").Value = 3575.6667
$ws.Range("J86").Value = 2158.8333
$ws.Range("K86").Value = 3575.6667
$ws.Range("L86").Value = 2158.8333
$ws.Range("M86").Value = -2452.6667
$ws.Range("N86").Value = -4404.8333
$ws.Range("H89").Value = 2867.25
$ws.Range("I89").Value = 3575.6667
$ws.Range("J89").Value = 2158.8333
$ws.Range("K89").Value = 17878.3335
$ws.Range("L89").Value = 10794.1665
$ws.Range("M89").Value = -12262.3335
$ws.Range("N89").Value = -22026.1665
$ws.Range("H105").Value = 3754.7778
$ws.Range("I105").Value = 3202.25
$ws.Range("K105").Value = 3202.25
$ws.Range("M105").Value = -1455.25
$ws.Range("H107").Value = 1460.4615
$ws.Range("I107").Value = 1221.5555
$ws.Range("K107").Value = 1221.5555
$ws.Range("M107").Value = 698.4445000000001
$ws.Range("H132").Value = 63019.668
$ws.Range("J132").Value = 63019.668
$ws.Range("L132").Value = 63019.668
$ws.Range("N132").Value = -73139.66800000001
$ws.Range("H134").Value = 1895
$ws.Range("I134").Value = 1895
$ws.Range("K134").Value = 5685
$ws.Range("M134").Value = -3150
$ws.Range("H140").Value = 67500
$ws.Range("J140").Value = 67500
$ws.Range("L140").Value = 67500
$ws.Range("N140").Value = -77860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2611.1904
$ws.Range("I31").Value = 1968.8572
$ws.Range("K31").Value = 1968.8572
$ws.Range("M31").Value = -1673.8572
$ws.Range("H34").Value = 2611.1904
$ws.Range("I34").Value = 1968.8572
$ws.Range("K34").Value = 1968.8572
$ws.Range("M34").Value = -1766.8572
$ws.Range("H58").Value = 3549.6
$ws.Range("I58").Value = 2899.5
$ws.Range("J58").Value = 3983
$ws.Range("K58").Value = 2899.5
$ws.Range("L58").Value = 3983
$ws.Range("M58").Value = -2696.5
$ws.Range("N58").Value = -4389
$ws.Range("H60").Value = 21974.125
$ws.Range("J60").Value = 21827.857
$ws.Range("L60").Value = 21827.857
$ws.Range("N60").Value = -22849.857
$ws.Range("H68").Value = 21950
$ws.Range("I68").Value = 4000
$ws.Range("J68").Value = 39900
$ws.Range("K68").Value = 4000
$ws.Range("L68").Value = 39900
$ws.Range("M68").Value = -3251
$ws.Range("N68").Value = -41398
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H71").Value = 21950
$ws.Range("I71").Value = 4000
$ws.Range("J71").Value = 39900
$ws.Range("K71").Value = 12000
$ws.Range("L71").Value = 119700
$ws.Range("M71").Value = -8256
$ws.Range("N71").Value = -127188
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H95").Value = 48666
$ws.Range("J95").Value = 48666
$ws.Range("L95").Value = 48666
$ws.Range("N95").Value = -54158
$ws.Range("H99").Value = 31119.7
$ws.Range("I99").Value = 6649.5
$ws.Range("J99").Value = 67825
$ws.Range("K99").Value = 6649.5
$ws.Range("L99").Value = 67825
$ws.Range("M99").Value = -5151.5
$ws.Range("N99").Value = -70821
$ws.Range("H126").Value = 31119.7
$ws.Range("I126").Value = 6649.5
$ws.Range("J126").Value = 67825
$ws.Range("K126").Value = 19948.5
$ws.Range("L126").Value = 203475
$ws.Range("M126").Value = -17478.5
$ws.Range("N126").Value = -208415
$ws.Range("H136").Value = 3549.6
$ws.Range("I136").Value = 2899.5
$ws.Range("J136").Value = 3983
$ws.Range("K136").Value = 8698.5
$ws.Range("L136").Value = 11949
$ws.Range("M136").Value = -6148.5
$ws.Range("N136").Value = -17049

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1250220
$ws.Range("J2").Value = 295
$ws.Range("L2").Value = 1770
$ws.Range("N2").Value = -1996
$ws.Range("H11").Value = 33682.332
$ws.Range("I11").Value = 50224.5
$ws.Range("K11").Value = 150673.5
$ws.Range("M11").Value = -150533.5
$ws.Range("H12").Value = 231.45833
$ws.Range("I12").Value = 315.36365
$ws.Range("J12").Value = 160.46153
$ws.Range("K12").Value = 946.09095
$ws.Range("L12").Value = 481.38459
$ws.Range("M12").Value = -773.09095
$ws.Range("N12").Value = -827.38459
$ws.Range("H17").Value = 1889.9
$ws.Range("J17").Value = 2428.5715
$ws.Range("L17").Value = 7285.7145
$ws.Range("N17").Value = -7623.7145
$ws.Range("H55").Value = 2422.5386
$ws.Range("I55").Value = 1799.75
$ws.Range("K55").Value = 5399.25
$ws.Range("M55").Value = -5222.25
$ws.Range("H98").Value = 631.3333
$ws.Range("J98").Value = 500
$ws.Range("L98").Value = 1500
$ws.Range("N98").Value = -4496
$ws.Range("H104").Value = 2000
$ws.Range("I104").Value = 2000
$ws.Range("K104").Value = 6000
$ws.Range("M104").Value = -3379
$ws.Range("H107").Value = 996.2105
$ws.Range("I107").Value = 671.6
$ws.Range("J107").Value = 1112.1428
$ws.Range("K107").Value = 2014.8
$ws.Range("L107").Value = 3336.4284
$ws.Range("M107").Value = -94.80000000000018
$ws.Range("N107").Value = -7176.428400000001
$ws.Range("H108").Value = 1298
$ws.Range("I108").Value = 1298
$ws.Range("K108").Value = 3894
$ws.Range("M108").Value = -1014
$ws.Range("H110").Value = 300
$ws.Range("I110").Value = 300
$ws.Range("K110").Value = 900
$ws.Range("M110").Value = 3190
$ws.Range("H132").Value = 6878.3335
$ws.Range("J132").Value = 8056
$ws.Range("L132").Value = 72504
$ws.Range("N132").Value = -77564
$ws.Range("H136").Value = 5705.8
$ws.Range("I136").Value = 4151.143
$ws.Range("K136").Value = 12453.429
$ws.Range("M136").Value = -7353.429
$ws.Range("H137").Value = 3497.8
$ws.Range("J137").Value = 3749.75
$ws.Range("L137").Value = 11249.25
$ws.Range("N137").Value = -21449.25
$ws.Range("H140").Value = 1669.1578
$ws.Range("I140").Value = 1669.1578
$ws.Range("K140").Value = 5007.4734
$ws.Range("M140").Value = 172.5266000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 12527.25
$ws.Range("I57").Value = 12527.25
$ws.Range("K57").Value = 12527.25
$ws.Range("M57").Value = -11707.25
$ws.Range("H64").Value = 69580.664
$ws.Range("J64").Value = 69580.664
$ws.Range("L64").Value = 69580.664
$ws.Range("N64").Value = -70076.664
$ws.Range("H67").Value = 69580.664
$ws.Range("J67").Value = 69580.664
$ws.Range("L67").Value = 69580.664
$ws.Range("N67").Value = -71296.664
$ws.Range("H70").Value = 8414
$ws.Range("I70").Value = 5492.25
$ws.Range("J70").Value = 10083.571
$ws.Range("K70").Value = 5492.25
$ws.Range("L70").Value = 10083.571
$ws.Range("M70").Value = -5222.25
$ws.Range("N70").Value = -10623.571
$ws.Range("H73").Value = 8414
$ws.Range("I73").Value = 5492.25
$ws.Range("J73").Value = 10083.571
$ws.Range("K73").Value = 5492.25
$ws.Range("L73").Value = 10083.571
$ws.Range("M73").Value = -4556.25
$ws.Range("N73").Value = -11955.571
$ws.Range("H95").Value = 46722.5
$ws.Range("I95").Value = 25000
$ws.Range("J95").Value = 53963.332
$ws.Range("K95").Value = 25000
$ws.Range("L95").Value = 53963.332
$ws.Range("M95").Value = -22254
$ws.Range("N95").Value = -59455.332
$ws.Range("H102").Value = 3661.3809
$ws.Range("I102").Value = 3270.25
$ws.Range("K102").Value = 3270.25
$ws.Range("M102").Value = -1648.25
$ws.Range("H113").Value = 2662.4
$ws.Range("I113").Value = 1933.3334
$ws.Range("J113").Value = 3756
$ws.Range("K113").Value = 1933.3334
$ws.Range("L113").Value = 3756
$ws.Range("M113").Value = 236.6666
$ws.Range("N113").Value = -8096
$ws.Range("H132").Value = 3747.8
$ws.Range("I132").Value = 3384.3635
$ws.Range("J132").Value = 4747.25
$ws.Range("K132").Value = 10153.0905
$ws.Range("L132").Value = 14241.75
$ws.Range("M132").Value = -7623.0905
$ws.Range("N132").Value = -19301.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 10000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 10000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 10000
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -10224
$ws.Range("H22").Value = 1832
$ws.Range("I22").Value = 849
$ws.Range("J22").Value = 2028.6
$ws.Range("K22").Value = 849
$ws.Range("L22").Value = 2028.6
$ws.Range("M22").Value = -554
$ws.Range("N22").Value = -2618.6
$ws.Range("H27").Value = 1832
$ws.Range("I27").Value = 849
$ws.Range("J27").Value = 2028.6
$ws.Range("K27").Value = 849
$ws.Range("L27").Value = 2028.6
$ws.Range("M27").Value = -742
$ws.Range("N27").Value = -2242.6
$ws.Range("H46").Value = 2444.762
$ws.Range("I46").Value = 820
$ws.Range("J46").Value = 2952.5
$ws.Range("K46").Value = 820
$ws.Range("L46").Value = 2952.5
$ws.Range("M46").Value = -632
$ws.Range("N46").Value = -3328.5
$ws.Range("H55").Value = 450.5238
$ws.Range("I55").Value = 342.58823
$ws.Range("K55").Value = 342.58823
$ws.Range("M55").Value = -169.58823
$ws.Range("H56").Value = 59249.75
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 59249.75
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 59249.75
$ws.Range("M56").ClearContents()
$ws.Range("N56").Value = -60631.75
$ws.Range("H68").Value = 2433.9443
$ws.Range("I68").Value = 2500.6365
$ws.Range("J68").Value = 2329.1428
$ws.Range("K68").Value = 2500.6365
$ws.Range("L68").Value = 2329.1428
$ws.Range("M68").Value = -1751.6365
$ws.Range("N68").Value = -3827.1428
$ws.Range("H71").Value = 2433.9443
$ws.Range("I71").Value = 2500.6365
$ws.Range("J71").Value = 2329.1428
$ws.Range("K71").Value = 12503.1825
$ws.Range("L71").Value = 11645.714
$ws.Range("M71").Value = -8759.182500000001
$ws.Range("N71").Value = -19133.714
$ws.Range("H82").Value = 2145.077
$ws.Range("I82").Value = 2245.4546
$ws.Range("K82").Value = 2245.4546
$ws.Range("M82").Value = -1884.4546
$ws.Range("H85").Value = 2145.077
$ws.Range("I85").Value = 2245.4546
$ws.Range("K85").Value = 2245.4546
$ws.Range("M85").Value = -997.4546
$ws.Range("H105").Value = 33807.5
$ws.Range("J105").Value = 33807.5
$ws.Range("L105").Value = 33807.5
$ws.Range("N105").Value = -40795.5
$ws.Range("H132").Value = 3939.4546
$ws.Range("I132").Value = 3699.3333
$ws.Range("J132").Value = 4227.6
$ws.Range("K132").Value = 11097.9999
$ws.Range("L132").Value = 12682.8
$ws.Range("M132").Value = -8567.999899999999
$ws.Range("N132").Value = -17742.8
$ws.Range("H136").Value = 3553.889
$ws.Range("J136").Value = 3399
$ws.Range("L136").Value = 10197
$ws.Range("N136").Value = -15297

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5947.4614
$ws.Range("I62").Value = 3981.8
$ws.Range("K62").Value = 3981.8
$ws.Range("M62").Value = -3357.8
$ws.Range("H65").Value = 5947.4614
$ws.Range("I65").Value = 3981.8
$ws.Range("K65").Value = 19909
$ws.Range("M65").Value = -16789
$ws.Range("H69").Value = 41140
$ws.Range("J69").Value = 41140
$ws.Range("L69").Value = 41140
$ws.Range("N69").Value = -42638
$ws.Range("H72").Value = 41140
$ws.Range("J72").Value = 41140
$ws.Range("L72").Value = 123420
$ws.Range("N72").Value = -130908
$ws.Range("H81").Value = 2224.75
$ws.Range("I81").Value = 1459.8
$ws.Range("J81").Value = 3499.6667
$ws.Range("K81").Value = 2919.6
$ws.Range("L81").Value = 6999.3334
$ws.Range("M81").Value = -1858.6
$ws.Range("N81").Value = -9121.3334
$ws.Range("H84").Value = 2224.75
$ws.Range("I84").Value = 1459.8
$ws.Range("J84").Value = 3499.6667
$ws.Range("K84").Value = 14598
$ws.Range("L84").Value = 34996.667
$ws.Range("M84").Value = -9294
$ws.Range("N84").Value = -45604.667
$ws.Range("H100").Value = 1784.9286
$ws.Range("I100").Value = 1998.4445
$ws.Range("J100").Value = 1400.6
$ws.Range("K100").Value = 3996.889
$ws.Range("L100").Value = 2801.2
$ws.Range("M100").Value = -3455.889
$ws.Range("N100").Value = -3883.2
$ws.Range("H112").Value = 43276.832
$ws.Range("J112").Value = 43276.832
$ws.Range("L112").Value = 43276.832
$ws.Range("N112").Value = -46230.832
$ws.Range("H123").Value = 67200
$ws.Range("J123").Value = 67200
$ws.Range("L123").Value = 67200
$ws.Range("N123").Value = -77000
$ws.Range("H132").Value = 4516.857
$ws.Range("I132").Value = 4492.7
$ws.Range("K132").Value = 13478.1
$ws.Range("M132").Value = -10948.1
$ws.Range("H136").Value = 4997.1665
$ws.Range("I136").Value = 3059.3125
$ws.Range("K136").Value = 9177.9375
$ws.Range("M136").Value = -6627.9375

